# Refresh the cryptocurrency price/volume snapshot (Price column D, Volume(1h) column E).
# Cells keep their original Text storage (matches the existing inline-string cells) by
# temporarily forcing a text number-format while writing, then restoring the cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "29.997.17"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -1.03%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.900.15"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -1.65%  "

$ws.Range("E4").Value = "  +0.00%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7418"
$dCell.Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "241.14"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3063"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -3.34%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "25.76"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -6.79%  "

$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06900"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -3.03%  "

$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08013"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -0.58%  "

$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7583"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -2.64%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "1.899.49"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -1.22%  "

$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "5.234"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -3.04%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "91.16"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -2.06%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "6.193"
$dCell.Style = "Normal"
$ws.Range("E16").Value = "  +2.84%  "

$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "30.002.19"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "14.00"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -3.79%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000007744"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  -2.46%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "237.25"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -5.80%  "

$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "1.001"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "2.146.22"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  -1.35%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "7.067"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  +5.79%  "

$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "9.312"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "166.82"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "

$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "18.82"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -1.35%  "

$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1262"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -2.74%  "

$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "2.035"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -6.85%  "

$ws.Range("E30").Value = "  -1.25%  "

$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "1.530"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -1.97%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "4.304"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -2.51%  "

$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "4.044"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -2.38%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05304"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  +1.37%  "

$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.288"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -2.47%  "

$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7396"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  -2.30%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "2.728"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -1.94%  "

$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.01934"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -0.68%  "

$ws.Range("E39").Value = "  -0.64%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "6.273"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -3.65%  "

$ws.Range("E41").Value = "  -1.62%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "72.96"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -6.62%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "1.957"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -1.05%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "0.8328"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -1.05%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "7.651"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "

$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = "101.16"
$dCell.Style = "Normal"
$ws.Range("E47").Value = "  -0.37%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "9.807"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -1.56%  "

$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "2.047.86"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "

$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "36.53"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -3.75%  "

$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1168"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -5.18%  "

